# Weekly CompStat (17th Precinct) report refresh: new crime-data week, volume/date
# header bump, and refreshed Week-to-Date/28-Day/YTD/2-Year figures + %Chg columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: volume number 48 -> 49, reporting week 11/25-12/1 -> 12/2-12/8 ---
$ws.Range("A8").Value = "Volume 31   Number  49"
$ws.Range("C9").Value = "Report Covering the Week  12/2/2024  Through  12/8/2024"

# --- Cells whose representation flips between a numeric figure and the text
#     placeholders "0" / "***.*" (used when a count or %Chg is not meaningful).
#     Row 14 (Murder) is untouched by this refresh, so we reuse its cells as
#     templates: copying them brings along the correct style + (for the text
#     placeholders) the correct shared string, matching how Excel itself keeps
#     formatting in sync when a cell's content type changes.
$ws.Range("C14").Copy($ws.Range("D16"))
$ws.Range("E14").Copy($ws.Range("E16"))
$ws.Range("J14").Copy($ws.Range("C17"))
$ws.Range("C17").Value = 1
$ws.Range("F14").Copy($ws.Range("D17"))
$ws.Range("D17").Value = 2
$ws.Range("K14").Copy($ws.Range("E17"))
$ws.Range("E17").Value = -50
$ws.Range("F14").Copy($ws.Range("C22"))
$ws.Range("C22").Value = 2
$ws.Range("C14").Copy($ws.Range("C28"))
$ws.Range("D14").Copy($ws.Range("D31"))
$ws.Range("H14").Copy($ws.Range("E31"))

# --- Remaining like-for-like numeric updates (counts + %Chg) ---
$ws.Range("L15").Value = -33.333333333333
$ws.Range("C16").Value = 2
$ws.Range("F16").Value = 4
$ws.Range("G16").Value = 7
$ws.Range("H16").Value = -42.857142857142
$ws.Range("I16").Value = 70
$ws.Range("K16").Value = -35.185185185185
$ws.Range("L16").Value = -19.540229885057
$ws.Range("M16").Value = 25
$ws.Range("N16").Value = -88.054607508532
$ws.Range("F17").Value = 9
$ws.Range("G17").Value = 4
$ws.Range("H17").Value = 125
$ws.Range("I17").Value = 107
$ws.Range("J17").Value = 100
$ws.Range("K17").Value = 7
$ws.Range("L17").Value = -0.925925925925
$ws.Range("M17").Value = 109.803921568627
$ws.Range("N17").Value = -30.519480519480
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 5
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = -58.333333333333
$ws.Range("I18").Value = 88
$ws.Range("J18").Value = 119
$ws.Range("K18").Value = -26.050420168067
$ws.Range("L18").Value = -49.425287356321
$ws.Range("M18").Value = 0
$ws.Range("N18").Value = -92.354474370112
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 19
$ws.Range("E19").Value = -47.368421052631
$ws.Range("F19").Value = 35
$ws.Range("G19").Value = 53
$ws.Range("H19").Value = -33.962264150943
$ws.Range("I19").Value = 550
$ws.Range("J19").Value = 653
$ws.Range("K19").Value = -15.773353751914
$ws.Range("L19").Value = -12.974683544303
$ws.Range("M19").Value = -20.289855072463
$ws.Range("N19").Value = -74.454249883883
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = -83.333333333333
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = -77.777777777777
$ws.Range("I20").Value = 25
$ws.Range("J20").Value = 57
$ws.Range("K20").Value = -56.140350877193
$ws.Range("L20").Value = -59.677419354838
$ws.Range("M20").Value = -16.666666666666
$ws.Range("N20").Value = -95.847176079734
$ws.Range("C21").Value = 16
$ws.Range("D21").Value = 29
$ws.Range("E21").Value = -44.827586206896
$ws.Range("G21").Value = 85
$ws.Range("H21").Value = -31.764705882352
$ws.Range("I21").Value = 854
$ws.Range("J21").Value = 1044
$ws.Range("K21").Value = -18.199233716475
$ws.Range("L21").Value = -20.779220779220
$ws.Range("M21").Value = -7.375271149674
$ws.Range("N21").Value = -81.669886241682
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 200
$ws.Range("I22").Value = 36
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 20
$ws.Range("M22").Value = 9.090909090909
$ws.Range("C24").Value = 17
$ws.Range("D24").Value = 18
$ws.Range("E24").Value = -5.555555555555
$ws.Range("F24").Value = 66
$ws.Range("G24").Value = 83
$ws.Range("H24").Value = -20.481927710843
$ws.Range("I24").Value = 875
$ws.Range("J24").Value = 1003
$ws.Range("K24").Value = -12.761714855433
$ws.Range("L24").Value = -25.784563189143
$ws.Range("M24").Value = 47.804054054054
$ws.Range("C25").Value = 13
$ws.Range("D25").Value = 14
$ws.Range("E25").Value = -7.142857142857
$ws.Range("F25").Value = 45
$ws.Range("G25").Value = 65
$ws.Range("H25").Value = -30.769230769230
$ws.Range("I25").Value = 673
$ws.Range("J25").Value = 824
$ws.Range("K25").Value = -18.325242718446
$ws.Range("L25").Value = -31.396534148827
$ws.Range("C26").Value = 4
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 14
$ws.Range("G26").Value = 19
$ws.Range("H26").Value = -26.315789473684
$ws.Range("I26").Value = 215
$ws.Range("J26").Value = 234
$ws.Range("K26").Value = -8.119658119658
$ws.Range("L26").Value = -3.587443946188
$ws.Range("M26").Value = -1.826484018264
$ws.Range("L27").Value = -29.411764705882
$ws.Range("L28").Value = -3.773584905660
$ws.Range("L33").Value = -50
